$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Kecepatan" (speed) criterion is being removed entirely, both its
# raw/rank column (C) and its normalized/value column (F). Delete the
# rightmost one first so column letters for the other doesn't shift.
$ws.Range("F1").EntireColumn.Delete()
$ws.Range("C1").EntireColumn.Delete()

# Two alternatives (A6 and A7, originally on rows 7 and 8) are removed.
$ws.Range("A7:A8").EntireRow.Delete()

# Recalculated SAW normalization / score values for the remaining
# alternatives, after the criterion removal and recalculation.
# Row 2 = alternative A1
$ws.Range("B2").Value = 0.2666666666666667
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0.3333333333333334
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.1333333333333334
$ws.Range("H2").Value = 0.4666666666666668

# Row 5 = alternative A4
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0.2173913043478261
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.07246376811594206
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0.5000000000000001
$ws.Range("H5").Value = 0.5724637681159421
